# Fruta / hortaliza, semanal
# Insert a new weekly block of 3 rows (Larga vida: Primera/Segunda/Tercera)
# just above the current row 308, shifting the rest of the table down by 3
# rows, then populate the new rows with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 308:310 - everything currently at 308 and below
# (through 399) shifts down to 311..402.
$ws.Rows("308:310").Insert()

# Common (unchanged across the block) field values, copied from the
# surrounding rows of this same market/product series.
$mercadoId   = 2
$mercado     = "Comercializadora del Agro de Limarí"
$region      = "Coquimbo"
$codigo      = 4
$productoId  = 100112020
$producto    = "Tomate"
$unidad      = '$/bandeja 18 kilos'
$provincia   = "Provincia de Limarí"
$cantidad    = 18
$tipo        = "Hortaliza"

# New row 308: Larga vida / Primera
$ws.Range("A308").Value2 = $mercadoId
$ws.Range("B308").Value2 = $mercado
$ws.Range("C308").Value2 = $region
$ws.Range("D308").Value2 = 44524
$ws.Range("E308").Value2 = $codigo
$ws.Range("F308").Value2 = $productoId
$ws.Range("G308").Value2 = $producto
$ws.Range("H308").Value2 = "Larga vida"
$ws.Range("I308").Value2 = "Primera"
$ws.Range("J308").Value2 = 2000
$ws.Range("K308").Value2 = 9000
$ws.Range("L308").Value2 = 10000
$ws.Range("M308").Value2 = 9500
$ws.Range("N308").Value2 = $unidad
$ws.Range("O308").Value2 = $provincia
$ws.Range("P308").Value2 = 528
$ws.Range("Q308").Value2 = $cantidad
$ws.Range("R308").Value2 = $tipo

# New row 309: Larga vida / Segunda
$ws.Range("A309").Value2 = $mercadoId
$ws.Range("B309").Value2 = $mercado
$ws.Range("C309").Value2 = $region
$ws.Range("D309").Value2 = 44524
$ws.Range("E309").Value2 = $codigo
$ws.Range("F309").Value2 = $productoId
$ws.Range("G309").Value2 = $producto
$ws.Range("H309").Value2 = "Larga vida"
$ws.Range("I309").Value2 = "Segunda"
$ws.Range("J309").Value2 = 1800
$ws.Range("K309").Value2 = 7000
$ws.Range("L309").Value2 = 8000
$ws.Range("M309").Value2 = 7500
$ws.Range("N309").Value2 = $unidad
$ws.Range("O309").Value2 = $provincia
$ws.Range("P309").Value2 = 417
$ws.Range("Q309").Value2 = $cantidad
$ws.Range("R309").Value2 = $tipo

# New row 310: Larga vida / Tercera
$ws.Range("A310").Value2 = $mercadoId
$ws.Range("B310").Value2 = $mercado
$ws.Range("C310").Value2 = $region
$ws.Range("D310").Value2 = 44524
$ws.Range("E310").Value2 = $codigo
$ws.Range("F310").Value2 = $productoId
$ws.Range("G310").Value2 = $producto
$ws.Range("H310").Value2 = "Larga vida"
$ws.Range("I310").Value2 = "Tercera"
$ws.Range("J310").Value2 = 1000
$ws.Range("K310").Value2 = 5000
$ws.Range("L310").Value2 = 6000
$ws.Range("M310").Value2 = 5500
$ws.Range("N310").Value2 = $unidad
$ws.Range("O310").Value2 = $provincia
$ws.Range("P310").Value2 = 306
$ws.Range("Q310").Value2 = $cantidad
$ws.Range("R310").Value2 = $tipo
